$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Insert a new row at 30 (pushes "sugar" from 30->31 and "tomato" from 31->32,
# and everything else shifts down by one row as well).
$ws.Rows(30).Insert()

# Populate the new row 30 with the "strawberries" ingredient data.
$ws.Range("A30").Value = "strawberries"
$ws.Range("B30").Value = 36
$ws.Range("C30").Value = 0.22
$ws.Range("D30").Value = 7.96
$ws.Range("E30").Value = 4.8600000000000003
$ws.Range("F30").Value = 0.64
$ws.Range("M30").Value = 3.7
$ws.Range("P30").Value = 59.6
$ws.Range("T30").Value = 161
$ws.Range("U30").Value = 17
$ws.Range("V30").Value = 23
$ws.Range("W30").Value = 12.5
$ws.Range("X30").Value = 0.26
$ws.Range("Y30").Value = 0.11
$ws.Range("Z30").Value = 0.035
$ws.Range("AA30").Value = 0.36799999999999999
$ws.Range("AC30").Value = 2

# Update the conditional-formatting ranges for columns G, H, I, J so they
# keep pointing at the same logical rows now that a row has been inserted.
$cf = $ws.Range("G1").FormatConditions
for ($i = 1; $i -le $cf.Count; $i++) {
    $fc = $cf.Item($i)
    if ($fc.AppliesTo.Address() -eq "`$G`$32:`$G`$1048576") {
        $fc.ModifyAppliesToRange($ws.Range("G33:G1048576"))
    }
}

$cf = $ws.Range("H1").FormatConditions
for ($i = 1; $i -le $cf.Count; $i++) {
    $fc = $cf.Item($i)
    if ($fc.AppliesTo.Address() -eq "`$H`$32:`$H`$1048576") {
        $fc.ModifyAppliesToRange($ws.Range("H33:H1048576"))
    }
}

$cf = $ws.Range("I1").FormatConditions
for ($i = 1; $i -le $cf.Count; $i++) {
    $fc = $cf.Item($i)
    if ($fc.AppliesTo.Address() -eq "`$I`$32:`$I`$1048576") {
        $fc.ModifyAppliesToRange($ws.Range("I33:I1048576"))
    }
}

$cf = $ws.Range("J1").FormatConditions
for ($i = 1; $i -le $cf.Count; $i++) {
    $fc = $cf.Item($i)
    if ($fc.AppliesTo.Address() -eq "`$J`$32:`$J`$1048576") {
        $fc.ModifyAppliesToRange($ws.Range("J33:J1048576"))
    }
}

# Move the sheet selection to match the recorded cursor position.
[void]$ws.Range("AB30").Select()
